$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from loginPage -> Sheet1
$ws.Name = "Sheet1"

# Remove the existing hyperlinks (Password@1 mailto links) and wipe the old data
$ws.Hyperlinks.Delete()
$ws.Cells.Delete()

# Rebuild the data grid with the Base_page / Auto_Constant rows
$ws.Range("A1").Value = "KAKBR0003"
$ws.Range("B1").Value = "Password@1"

$ws.Range("A2").Value = "objName"
$ws.Range("B2").Value = "objPath"

$ws.Range("A3").Value = "login"
$ws.Range("B3").Value = "//a[@class='loginTop']"

$ws.Range("A4").Value = "user_name"
$ws.Range("B4").Value = "UserName"

$ws.Range("A5").Value = "user_pwd"
$ws.Range("B5").Value = "Password"

$ws.Range("A6").Value = "loginBtn"
$ws.Range("B6").Value = "//span[normalize-space()='Login']"

# Style the Password cell with a plain (non-hyperlink) dark-grey font
# (no underline, no theme scheme, explicit RGB color FF202124)
$f = $ws.Range("B5").Font
$f.Color = 2367776
$f.Name = "Calibri"

# Widen column B (33 chars, accounting for the engine's internal padding offset) and update the view/selection
$ws.Columns("B").ColumnWidth = 32 + 1/6
$ws.Range("C19").Select()

# Resize the workbook window (best effort)
$win = $excel.Windows.Item(1)
$win.Width = 15675
$win.Height = 6615
